{"js": "// The edit replaces the date line and 25 of the division-fact answers\n// inside the worksheet table (5 filled rows x 5 columns, with blank\n// spacer rows in between). We walk the paragraphs of the document body\n// in document order and rewrite every non-empty one according to the\n// positional mapping derived from the OOXML diff (the Nth non-blank\n// paragraph in the \"before\" doc maps to the Nth replacement value).\n\nconst replacements = [\n  \"2024-06-21 Friday\",\n  \"41\u00f78=5, 1\",\n  \"81\u00f77=11, 4\",\n  \"69\u00f76=11, 3\",\n  \"26\u00f72=13, 0\",\n  \"77\u00f78=9, 5\",\n  \"24\u00f74=6, 0\",\n  \"72\u00f73=24, 0\",\n  \"77\u00f72=38, 1\",\n  \"95\u00f78=11, 7\",\n  \"99\u00f74=24, 3\",\n  \"10\u00f79=1, 1\",\n  \"75\u00f79=8, 3\",\n  \"54\u00f76=9, 0\",\n  \"47\u00f72=23, 1\",\n  \"31\u00f77=4, 3\",\n  \"41\u00f78=5, 1\",\n  \"84\u00f76=14, 0\",\n  \"63\u00f74=15, 3\",\n  \"70\u00f73=23, 1\",\n  \"12\u00f76=2, 0\",\n  \"96\u00f72=48, 0\",\n  \"82\u00f78=10, 2\",\n  \"65\u00f74=16, 1\",\n  \"14\u00f74=3, 2\",\n  \"63\u00f73=21, 0\",\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet replaceIdx = 0;\nfor (let i = 0; i < paragraphs.items.length && replaceIdx < replacements.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text.trim().length === 0) {\n    continue;\n  }\n  para.insertText(replacements[replaceIdx], \"Replace\");\n  replaceIdx++;\n}\n\nawait context.sync();\n", "ps1": "# The edit updates the date heading and the 25 division-fact answers\n# that live in the worksheet table (5 filled rows x 5 columns, with\n# blank spacer rows interleaved). We address each filled cell directly\n# by its (row, column) position in Tables(1) so duplicate answer\n# strings (e.g. \"18\u00f79=2, 0\" appearing twice, each mapping to a\n# different replacement) are handled unambiguously.\n\n$d = $word.ActiveDocument\n\n# Date heading (first paragraph, outside the table).\n$d.Paragraphs(1).Range.Text = \"2024-06-21 Friday\"\n\n$t = $d.Tables(1)\n\n$values = @(\n    @(\"41\u00f78=5, 1\", \"81\u00f77=11, 4\", \"69\u00f76=11, 3\", \"26\u00f72=13, 0\", \"77\u00f78=9, 5\"),\n    @(\"24\u00f74=6, 0\", \"72\u00f73=24, 0\", \"77\u00f72=38, 1\", \"95\u00f78=11, 7\", \"99\u00f74=24, 3\"),\n    @(\"10\u00f79=1, 1\", \"75\u00f79=8, 3\", \"54\u00f76=9, 0\", \"47\u00f72=23, 1\", \"31\u00f77=4, 3\"),\n    @(\"41\u00f78=5, 1\", \"84\u00f76=14, 0\", \"63\u00f74=15, 3\", \"70\u00f73=23, 1\", \"12\u00f76=2, 0\"),\n    @(\"96\u00f72=48, 0\", \"82\u00f78=10, 2\", \"65\u00f74=16, 1\", \"14\u00f74=3, 2\", \"63\u00f73=21, 0\")\n)\n\n$rows = @(1, 5, 9, 13, 17)\n\nfor ($i = 0; $i -lt $rows.Count; $i++) {\n    $rowIndex = $rows[$i]\n    $rowValues = $values[$i]\n    for ($c = 1; $c -le 5; $c++) {\n        $cell = $t.Cell($rowIndex, $c)\n        $cell.Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
